# Data Driving Invalid Login Test Script
#
# Adds a new "InvalidLogin" worksheet (after the existing "ValidLogin" sheet)
# that holds a bad username/password pair, used for negative-path testing.

$wb = $excel.ActiveWorkbook

# Add the new worksheet right after the last existing sheet (ValidLogin).
$newSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "InvalidLogin"

# Same header row as ValidLogin, plus the bogus credentials to drive the
# invalid-login test case.
$newSheet.Range("A1").Value = "UserName"
$newSheet.Range("B1").Value = "Password"
$newSheet.Range("A2").Value = "abc"
$newSheet.Range("B2").Value = "xyz"

# Make the new sheet the active/selected tab, with its own zoom level and
# cell selection, mirroring how it would look after being authored in Excel.
$newSheet.Activate()
$excel.ActiveWindow.Zoom = 220
$newSheet.Range("A3").Select()
